$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("N3").Value = 5
$ws.Range("C4").Value = 4
$ws.Range("K8").Value = 5
$ws.Range("K17").Value = 6
